$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 15:35"

# --- Update Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1387407
$ws.Range("C4").Value = 1573
$ws.Range("E4").Value = 1043273
$ws.Range("G4").Value = 114
$ws.Range("H4").Value = 81909

# --- Update Brasil (row 11) ---
$ws.Range("B11").Value = 169906
$ws.Range("C11").Value = 763
$ws.Range("E11").Value = 90827
$ws.Range("G11").Value = 70
$ws.Range("H11").Value = 11695

# --- Swap Singapur/Bielorrusia (rows 29/30) and update their stats ---
$ws.Range("A29").Value = "Bielorrusia"
$ws.Range("B29").Value = 24873
$ws.Range("C29").Value = 967
$ws.Range("D29").Value = 6974
$ws.Range("E29").Value = 17757
$ws.Range("F29").Value = 92
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 142

$ws.Range("A30").Value = "Singapur"
$ws.Range("B30").Value = 24671
$ws.Range("C30").Value = 884
$ws.Range("D30").Value = 3225
$ws.Range("E30").Value = 21425
$ws.Range("F30").Value = 24
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 21

# --- Swap San Bartolome/Sahara Occidental (rows 215/216); stats unchanged ---
$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("A216").Value = "San Bartolome"
